$p = $ppt.ActivePresentation

$CR = [char]13
$LF = [string][char]10
$LDQ = [char]8220   # left double quotation mark "
$RDQ = [char]8221   # right double quotation mark "
$ELL = [char]8230   # horizontal ellipsis …

# ------------------------------------------------------------------
# 1) Slide 9 notes (notesSlide7.xml): append two new paragraphs after
#    the existing note text - a blank paragraph, then the new remark.
# ------------------------------------------------------------------
$slide9 = $p.Slides.Item(9)
$notes9Shape = $slide9.NotesPage.Shapes.Item(2)
$notes9Range = $notes9Shape.TextFrame.TextRange
$notes9Parts = [System.Collections.ArrayList]@($notes9Range.Text.Split($CR))

$notes9NewPara = "It is IMPORTANT to highlight that the function does not complete due to " + $LDQ + "await" + $RDQ + " but the control IMMEDIATELY returns to the current thread " + $ELL + ". Which ALLOWS THE CURRRENT THREAD TO COMPLETE."
[void]$notes9Parts.Add("")
[void]$notes9Parts.Add($notes9NewPara)

$notes9Range.Text = $notes9Parts -join $LF

# ------------------------------------------------------------------
# 2) Slide 10 notes (notesSlide8.xml):
#    - replace paragraph 2 ("When one request starts, this code
#      proceeds to the next one.") with two expanded paragraphs plus
#      a trailing blank paragraph
#    - add a blank paragraph + a new paragraph after the
#      "Transmission is out-of-order..." paragraph
# ------------------------------------------------------------------
$slide10 = $p.Slides.Item(10)
$notes10Shape = $slide10.NotesPage.Shapes.Item(2)
$notes10Range = $notes10Shape.TextFrame.TextRange
$notes10Parts = [System.Collections.ArrayList]@($notes10Range.Text.Split($CR))

$notes10NewPara2 = "When you call makeOneGetRequest, one request starts, the function makes the axios.get request, suspends itself due to " + $LDQ + "await" + $RDQ + ", a promise is returned immediately which gives the control back to the calling function " + $ELL + "."
$notes10NewPara3 = "this code proceeds to the second makeOneGetRequest call."

$oldPara2Index = -1
for ($i = 0; $i -lt $notes10Parts.Count; $i++) {
    if ($notes10Parts[$i] -eq "When one request starts, this code proceeds to the next one.") {
        $oldPara2Index = $i
    }
}
if ($oldPara2Index -ge 0) {
    $notes10Parts[$oldPara2Index] = $notes10NewPara2
    $notes10Parts.Insert($oldPara2Index + 1, $notes10NewPara3)
    $notes10Parts.Insert($oldPara2Index + 2, "")
}

$transIndex = -1
for ($i = 0; $i -lt $notes10Parts.Count; $i++) {
    if ($notes10Parts[$i] -eq "Transmission is out-of-order:  Request 2 evidently reached the server before Request 1.") {
        $transIndex = $i
    }
}
$notes10NewExtraPara = "Three tasks were suspended " + $ELL + ". They get executed in some order"
if ($transIndex -ge 0) {
    $notes10Parts.Insert($transIndex + 1, "")
    $notes10Parts.Insert($transIndex + 2, $notes10NewExtraPara)
}

$notes10Range.Text = $notes10Parts -join $LF

# ------------------------------------------------------------------
# 3) Slide 9 body text (slide9.xml, "Content Placeholder 4"): change
#    " is resolved)." to " is resolved or rejected)." while keeping
#    the surrounding runs/formatting untouched.
# ------------------------------------------------------------------
$contentShape = $slide9.Shapes.Item(4)
$contentRange = $contentShape.TextFrame.TextRange
$fullText = $contentRange.Text
$target = " is resolved)."
$startPos = $fullText.IndexOf($target) + 1
$sub = $contentRange.Characters($startPos, $target.Length)
$sub.Text = " is resolved or rejected)."
